# Applies the "answers-of-two-digit-divison" update:
#  - bump the date line
#  - rewrite the division answers in rows 1, 5 (all 5 cells) and rows 12 (cells 2-5)
#    and the row that is row 17 before the row-block shuffle
#  - remove the block of rows that held "37÷9=4, 1 | 27÷3=9, 0 | 77÷3=25, 2 | 95÷4=23, 3 | 35÷2=17, 1"
#    (together with the 3 preceding blank rows)
#  - append a new block of 3 blank rows + one content row at the end of the table

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Date heading
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("2023-11-28 Tuesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2023-11-29 Wednesday", 2)

$t = $d.Tables.Item(1)

# ---------------------------------------------------------------------------
# Helper: replace the text of one cell while keeping its existing run/para
# formatting (the cell already carries rFonts=TimeNewRoman, sz=30, jc=left,
# so a plain text assignment keeps that formatting intact).
# ---------------------------------------------------------------------------
function Set-CellText($cell, $text) {
    $cell.Range.Text = $text
}

# Helper: populate a *blank* cell (no existing run) with properly formatted
# text by inserting a small WordprocessingML fragment - this is what lets us
# stamp the same rFonts/sz/jc combo onto newly-added rows.
function Set-BlankCellText($cell, $text) {
    $frag = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:jc w:val="left"/></w:pPr><w:r><w:rPr><w:rFonts w:ascii="TimeNewRoman" w:hAnsi="TimeNewRoman"/><w:sz w:val="30"/></w:rPr><w:t>' + $text + '</w:t></w:r></w:p>'
    $cell.Range.InsertXML($frag)
}

# ---------------------------------------------------------------------------
# 2. Row 1 (all five cells)
# ---------------------------------------------------------------------------
$row = $t.Rows.Item(1)
Set-CellText $row.Cells.Item(1) "93÷7=13, 2"
Set-CellText $row.Cells.Item(2) "30÷5=6, 0"
Set-CellText $row.Cells.Item(3) "19÷4=4, 3"
Set-CellText $row.Cells.Item(4) "18÷2=9, 0"
Set-CellText $row.Cells.Item(5) "51÷7=7, 2"

# ---------------------------------------------------------------------------
# 3. Row 5 (all five cells)
# ---------------------------------------------------------------------------
$row = $t.Rows.Item(5)
Set-CellText $row.Cells.Item(1) "72÷5=14, 2"
Set-CellText $row.Cells.Item(2) "99÷7=14, 1"
Set-CellText $row.Cells.Item(3) "94÷9=10, 4"
Set-CellText $row.Cells.Item(4) "71÷7=10, 1"
Set-CellText $row.Cells.Item(5) "75÷7=10, 5"

# ---------------------------------------------------------------------------
# 4. Drop the 3 blank rows (6-8) plus the content row (9) that held
#    37÷9=4,1 | 27÷3=9,0 | 77÷3=25,2 | 95÷4=23,3 | 35÷2=17,1
# ---------------------------------------------------------------------------
for ($k = 0; $k -lt 4; $k++) {
    $t.Rows.Item(6).Delete()
}

# ---------------------------------------------------------------------------
# 5. What used to be row 12 is now row 9; its first cell (47÷6=7, 5) is left
#    untouched, only cells 2-5 change.
# ---------------------------------------------------------------------------
$row = $t.Rows.Item(9)
Set-CellText $row.Cells.Item(2) "30÷8=3, 6"
Set-CellText $row.Cells.Item(3) "87÷8=10, 7"
Set-CellText $row.Cells.Item(4) "97÷5=19, 2"
Set-CellText $row.Cells.Item(5) "57÷3=19, 0"

# ---------------------------------------------------------------------------
# 6. What used to be row 16 is now row 13 (all five cells change)
# ---------------------------------------------------------------------------
$row = $t.Rows.Item(13)
Set-CellText $row.Cells.Item(1) "95÷5=19, 0"
Set-CellText $row.Cells.Item(2) "19÷6=3, 1"
Set-CellText $row.Cells.Item(3) "85÷8=10, 5"
Set-CellText $row.Cells.Item(4) "57÷2=28, 1"
Set-CellText $row.Cells.Item(5) "25÷9=2, 7"

# ---------------------------------------------------------------------------
# 7. Append a new block: 3 blank rows, then one content row with
#    55÷5=11,0 | 30÷6=5,0 | 85÷6=14,1 | 23÷3=7,2 | 15÷6=2,3
# ---------------------------------------------------------------------------
$t.Rows.Add() | Out-Null
$t.Rows.Add() | Out-Null
$t.Rows.Add() | Out-Null
$t.Rows.Add() | Out-Null

$newRow = $t.Rows.Item($t.Rows.Count)
Set-BlankCellText $newRow.Cells.Item(1) "55÷5=11, 0"
Set-BlankCellText $newRow.Cells.Item(2) "30÷6=5, 0"
Set-BlankCellText $newRow.Cells.Item(3) "85÷6=14, 1"
Set-BlankCellText $newRow.Cells.Item(4) "23÷3=7, 2"
Set-BlankCellText $newRow.Cells.Item(5) "15÷6=2, 3"
